# fix fastq filedate in jp 092120
# The libraryDate value "08.30.20" (shared by rows 19-26 in column A) is
# corrected to "08.27.20". Force a text number format before/after the
# assignment so Excel's auto date-detection doesn't turn the literal
# "08.27.20" string into a date serial number - the source file stores it
# as plain text.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$dateRange = $ws.Range("A19:A26")
$dateRange.NumberFormat = "@"
$dateRange.Value = "08.27.20"
$dateRange.NumberFormat = "General"

# Row heights for rows 20-26 were nudged down to 13.8pt.
$ws.Rows("20:26").RowHeight = 13.8

# Window/view state: zoomed out to 100% and scrolled so the edited rows
# are visible, leaving the selection on A26.
$excel.ActiveWindow.Zoom = 100
$ws.Range("A26").Select()
